$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old TOTAL column (L1:L2)
$ws.Range("L1:L2").Clear()

# Row 4: variance per activity (squared deviation)
$ws.Range("A4").Value = "Varianza"
$ws.Range("B4").Formula = "=B3^2"
$ws.Range("C4:K4").Formula = "=C3^2"

# Row 6: headers for route table
$ws.Range("A6").Value = "Rutas"
$ws.Range("B6").Value = "Tiempo"
$ws.Range("C6").Value = "Varianza"
$ws.Range("D6").Value = "Desviacion estandar"

# Row 7
$ws.Range("A7").Value = "A, E, G, I y J"
$ws.Range("B7").Formula = "=B2 + F2 + H2 + J2 + K2"

# Row 8
$ws.Range("A8").Value = "A, B, D, F, G y J"
$ws.Range("B8").Formula = " B2 + C2 + E2 + G2 + H2 + K2"

# Row 9
$ws.Range("A9").Value = "A, C, E, G, H y J"
$ws.Range("B9").Formula = " B2 + D2 + F2 + H2 + I2 + K2"
$ws.Range("C9").Formula = "=B4+D4+F4+H4+I4+K4"
$ws.Range("D9").Formula = "=SQRT(C9)"

# Row 10
$ws.Range("A10").Value = "A, C, D, E, F, H y J"
$ws.Range("B10").Formula = " B2 + D2 + E2 + F2 + G2 + I2 + K2"

# Row 11
$ws.Range("A11").Value = "TE ="
$ws.Range("B11").Formula = "=MAX(B7:B10)"
$ws.Range("C11").Value = "P = "
$ws.Range("D11").Value = 0.5
